$d = $word.ActiveDocument

# 1) Remove the entire "July 22nd, 2019" paragraph (including its paragraph mark).
foreach ($p in @($d.Paragraphs)) {
    if ($p.Range.Text -eq "July 22nd, 2019`r") {
        $p.Range.Delete()
        break
    }
}

# 2) Move the "_GoBack" bookmark from before "numbers, and the date..." to
#    wrap "stuff. " within the first body paragraph
#    ("...purchase stuff. In addition...").
if ($d.Bookmarks.Exists("_GoBack")) {
    [void]$d.Bookmarks.Item("_GoBack").Delete()
}

$rng = $d.Content
[void]$rng.Find.Execute("stuff. In addition", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$target = $d.Range($start, $start + 7)
[void]$d.Bookmarks.Add("_GoBack", $target)
